$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-25
# from serial value 46060 (2026-02-07) to 46061 (2026-02-08)
for ($row = 2; $row -le 25; $row++) {
    $cell = $ws.Range("C$row")
    if ($cell.Value2 -eq 46060) {
        $cell.Value2 = 46061
    }
}
